# Append the 2025-04-09 price observation to every Solar_Prices sheet.
# Each worksheet holds a two-column (Date, Price) table running through
# row 38 (2025-04-08). This adds a new row 39 dated 2025-04-09, carrying
# forward the same price as the prior day (matching the source data).

$wb = $excel.ActiveWorkbook

$newRow = @{
    "N-Dense"                  = "40"
    "N-Type"                   = "41.5"
    "N-type Wafer"              = "1.28"
    "Cell Topcon 183mm"         = "0.303"
    "Module Topcon 183mm"       = "0.1"
    "Silver Rear_side"          = "5,039"
    "Silver Busbar front-side"  = "7,543"
    "Silver finger front-side"  = "7,593"
    "USD_CNY"                   = "7.3569"
}

$sheetOrder = @(
    "N-Dense",
    "N-Type",
    "N-type Wafer",
    "Cell Topcon 183mm",
    "Module Topcon 183mm",
    "Silver Rear_side",
    "Silver Busbar front-side",
    "Silver finger front-side",
    "USD_CNY"
)

foreach ($name in $sheetOrder) {
    $ws = $wb.Worksheets.Item($name)

    $dateCell = $ws.Range("A39")
    $priceCell = $ws.Range("B39")

    # Force text storage (matching every other data row in the column)
    # instead of Excel's automatic date/number inference.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = "2025-04-09"
    $priceCell.Value = $newRow[$name]

    # Drop the temporary text formatting so the new cells end up with
    # the same (default) style as the rest of the table.
    $ws.Range("A39:B39").ClearFormats()
}
